# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect refreshed output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1166
$ws1.Range("F3").Value = 1346
$ws1.Range("F4").Value = 302
$ws1.Range("F5").Value = 1022
$ws1.Range("F6").Value = 10569
$ws1.Range("F11").Value = 686
$ws1.Range("F12").Value = 12006
$ws1.Range("F13").Value = 12415
$ws1.Range("F15").Value = 117

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1166
$ws4.Range("F4").Value = 1346
$ws4.Range("F5").Value = 302
$ws4.Range("F6").Value = 1022
$ws4.Range("F7").Value = 10569
$ws4.Range("F12").Value = 686
$ws4.Range("F13").Value = 12006
$ws4.Range("F14").Value = 12415
$ws4.Range("F16").Value = 117
